$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "submitter_email" column (H) is being removed from the sheet's data:
# the header and all data values in column H are cleared out, while the
# column itself (and its formatting) remains in place.
$ws.Range("H1:H4").ClearContents()

# Update the saved selection to match the edited file (H1:H10 selected,
# active cell H1).
$ws.Range("H1:H10").Select()
